# Update presentation for "Winter Session 2025":
#  - Refresh the cached "datetimeFigureOut" footer-date field text
#    (7/25/24 -> 1/17/25) across the slide master, every slide layout,
#    and the notes master.
#  - Rename the title-slide subtitle from "... SummerSession 2024" to
#    "... WinterSession 2025".

$p = $ppt.ActivePresentation
$newDate = "1/17/25"

# --- Slide master: "Date Placeholder" shape -------------------------------
foreach ($sh in $p.SlideMaster.Shapes) {
    if ($sh.Name -like "*Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = $newDate
    }
}

# --- Every slide layout: "Date Placeholder" shape --------------------------
foreach ($lay in $p.SlideMaster.CustomLayouts) {
    foreach ($sh in $lay.Shapes) {
        if ($sh.Name -like "*Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = $newDate
        }
    }
}

# --- Notes master date field ------------------------------------------------
# (Shapes(...).TextFrame.TextRange.Text is a no-op for the NotesMaster in
# this runtime, so go through the HeadersFooters / DateAndTime object
# instead, which does persist the edit.)
$p.NotesMaster.HeadersFooters.DateAndTime.Text = $newDate

# --- Slide 1 title subtitle: SummerSession 2024 -> WinterSession 2025 ------
$s1 = $p.Slides.Item(1)
foreach ($sh in $s1.Shapes) {
    if ($sh.HasTextFrame -and ($sh.TextFrame.TextRange.Text -like "*SummerSession*")) {
        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        $i1 = $full.IndexOf("SummerSession")
        $tr.Characters($i1 + 1, "SummerSession".Length).Text = "WinterSession"

        $full = $tr.Text
        $i2 = $full.IndexOf(" 2024")
        $tr.Characters($i2 + 1, " 2024".Length).Text = " 2025"
    }
}
